# ConfigurationDetails.xlsx - "Added system path for chromedriver"
#
# The Chrome Driver row used to store the contributor's own absolute
# machine path to chromedriver.exe. Replace it with a relative/system
# path, keep the row's formatting consistent with the other driver/
# environment rows (hyperlink style + auto row height), and register a
# hyperlink for it just like the other URL/path cells on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

$newChromePath = "\\exe\\chromedriver.exe"

# 1) Update the Chrome Driver path (row 17, column B).
$ws.Range("B17").Value2 = $newChromePath

# 2) The old value needed a wrapped 2-line row; the new short value
#    fits on one line, so let the row shrink back to the default height.
$ws.Rows(17).AutoFit()

# 3) Make the cell a hyperlink, matching the rest of the "Driver Path" /
#    "Web Portal Environment" rows above it.
$ws.Hyperlinks.Add($ws.Range("B17"), $newChromePath) | Out-Null

# 4) Hyperlinks.Add() stamps its own style on the cell - reapply the same
#    "Hyperlink" look already used by the sibling cells (e.g. B2) so B17
#    matches them exactly instead of getting an ad-hoc one.
$ws.Range("B2").Copy()
$ws.Range("B17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 5) Leave the selection on the row that was edited.
$ws.Range("C17").Select()
